$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.697.17"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.486.45"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.38"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.430"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.090.51"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.75"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.729.65"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.480.57"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.01"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "391.53"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.91"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.65"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.533"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("E25").Value = "  -2.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.14"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.18%  "
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.18"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.42"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.18%  "
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.62"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.30"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.75"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.64"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0737"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.11"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.794.74"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.99"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.64"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0301"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "335.75"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.26"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("E51").Value = "  -2.63%  "
